$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.482.97"
$ws.Range("D3").Value = "2.631.77"
$ws.Range("E3").Value = "  -1.65%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "594.72"
$ws.Range("E5").Value = "  -0.64%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "167.93"
$ws.Range("E6").Value = "  +0.95%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.533"
$ws.Range("E8").Value = "  -2.44%  "
$ws.Range("D9").Value = "2.631.26"
$ws.Range("E10").Value = "  -3.46%  "
$ws.Range("E11").Value = "  +1.30%  "
$ws.Range("E12").Value = "  +0.59%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.22"
$ws.Range("E13").Value = "  -0.19%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.66"
$ws.Range("E14").Value = "  -0.75%  "
$ws.Range("D15").Value = "3.110.97"
$ws.Range("E15").Value = "  -1.66%  "
$ws.Range("E16").Value = "  -2.16%  "
$ws.Range("D17").Value = "67.326.44"
$ws.Range("E17").Value = "  -0.53%  "
$ws.Range("D18").Value = "2.612.04"
$ws.Range("E18").Value = "  -2.36%  "
$ws.Range("E19").Value = "  +1.92%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.04"
$ws.Range("E20").Value = "  +3.87%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "357.20"
$ws.Range("E21").Value = "  -1.94%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.31"
$ws.Range("E22").Value = "  -1.89%  "
$ws.Range("E23").Value = "  -3.28%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.95"
$ws.Range("E24").Value = "  -4.57%  "
$ws.Range("E25").Value = "  -0.01%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.33"
$ws.Range("E26").Value = "  +2.10%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "69.67"
$ws.Range("E27").Value = "  -1.93%  "
$ws.Range("D28").Value = "2.767.97"
$ws.Range("E28").Value = "  -1.88%  "
$ws.Range("E29").Value = "  +0.40%  "
$ws.Range("E30").Value = "  -1.80%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "548.05"
$ws.Range("E31").Value = "  -1.70%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.94"
$ws.Range("E32").Value = "  -1.05%  "
$ws.Range("E33").Value = "  -3.10%  "
$ws.Range("E34").Value = "  -2.00%  "
$ws.Range("E35").Value = "  +4.40%  "
$ws.Range("E36").Value = "  +0.02%  "
$ws.Range("E37").Value = "  -3.92%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "157.01"
$ws.Range("E38").Value = "  +0.68%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "19.02"
$ws.Range("E39").Value = "  -2.67%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.364"
$ws.Range("E40").Value = "  -2.41%  "
$ws.Range("E41").Value = "  -0.59%  "
$ws.Range("E42").Value = "  +1.97%  "
$ws.Range("E43").Value = "  -1.75%  "
$ws.Range("E45").Value = "  -3.92%  "
$ws.Range("D46").Value = "0.0₆0298"
$ws.Range("E46").Value = "  -0.66%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "153.09"
$ws.Range("E47").Value = "  -0.32%  "
$ws.Range("E48").Value = "  -2.04%  "
$ws.Range("E49").Value = "  -1.48%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.70"
$ws.Range("E50").Value = "  -1.29%  "
$ws.Range("E51").Value = "  -1.17%  "
